# Atualização automática de CRUZ_ALTA.xlsx
$wb = $excel.ActiveWorkbook

# Excel requires at least one sheet to remain visible; disable the
# confirmation prompt for deleting a sheet (no-op in this headless runtime,
# but kept for fidelity with real Excel automation scripts).
$excel.DisplayAlerts = $false

# Remove the "Desarquivamentos Pendentes" worksheet entirely.
$wb.Worksheets.Item("Desarquivamentos Pendentes").Delete()

# Rename the remaining sheets to their new (upper-cased / accented) titles.
$wb.Worksheets.Item("Paineis DARQ").Name = "PAINEIS DARQ"
$wb.Worksheets.Item("Recolhimento x Eliminacao").Name = "RECOLHIMENTO X ELIMINAÇÃO"

$excel.DisplayAlerts = $true
